$wb = $excel.ActiveWorkbook

# ---- Caso1 ----
$ws = $wb.Worksheets.Item("Caso1")

# Headers for new AutoML columns (H1, I1, J1)
$ws.Cells.Item(1,8).Value = "Valori stimati autogloun"
$ws.Cells.Item(1,9).Value = "Valori stimati h2o"
$ws.Cells.Item(1,10).Value = "Valori stimati autosklearn"

# Copy formatting (no explicit style) from G1 onto the new header cells
$ws.Cells.Item(1,7).Copy()
$ws.Range($ws.Cells.Item(1,8), $ws.Cells.Item(1,10)).PasteSpecial(-4122)

# Update existing MLP column (G) and populate new AutoML columns (H,I,J) for rows 2-19
$ws.Cells.Item(2,7).Value = 1.009679436683655
$ws.Cells.Item(2,8).Value = 1.00510585308075
$ws.Cells.Item(2,9).Value = 1.004770611153234
$ws.Cells.Item(2,10).Value = 1.004751645028591
$ws.Cells.Item(3,7).Value = 1.006496548652649
$ws.Cells.Item(3,8).Value = 1.005297422409058
$ws.Cells.Item(3,9).Value = 1.004581883632716
$ws.Cells.Item(3,10).Value = 1.004619345068932
$ws.Cells.Item(4,7).Value = 1.007214665412903
$ws.Cells.Item(4,8).Value = 1.00562310218811
$ws.Cells.Item(4,9).Value = 1.004483717599561
$ws.Cells.Item(4,10).Value = 1.00461733341217
$ws.Cells.Item(5,7).Value = 1.005931735038757
$ws.Cells.Item(5,8).Value = 1.005711674690247
$ws.Cells.Item(5,9).Value = 1.004488741128931
$ws.Cells.Item(5,10).Value = 1.004608931019902
$ws.Cells.Item(6,7).Value = 1.006258487701416
$ws.Cells.Item(6,8).Value = 1.005583643913269
$ws.Cells.Item(6,9).Value = 1.004502521130789
$ws.Cells.Item(6,10).Value = 1.004726594313979
$ws.Cells.Item(7,7).Value = 1.005318164825439
$ws.Cells.Item(7,8).Value = 1.00560200214386
$ws.Cells.Item(7,9).Value = 1.004509439153644
$ws.Cells.Item(7,10).Value = 1.004624389111996
$ws.Cells.Item(8,7).Value = 1.005589842796326
$ws.Cells.Item(8,8).Value = 1.005555033683777
$ws.Cells.Item(8,9).Value = 1.004532997790389
$ws.Cells.Item(8,10).Value = 1.004691377282143
$ws.Cells.Item(9,7).Value = 1.006085872650146
$ws.Cells.Item(9,8).Value = 1.005752563476562
$ws.Cells.Item(9,9).Value = 1.004500216475008
$ws.Cells.Item(9,10).Value = 1.004654239863157
$ws.Cells.Item(10,7).Value = 1.006041049957275
$ws.Cells.Item(10,8).Value = 1.005478382110596
$ws.Cells.Item(10,9).Value = 1.004495110606235
$ws.Cells.Item(10,10).Value = 1.004611391574144
$ws.Cells.Item(11,7).Value = 1.005326509475708
$ws.Cells.Item(11,8).Value = 1.005604267120361
$ws.Cells.Item(11,9).Value = 1.004498267014363
$ws.Cells.Item(11,10).Value = 1.004671387374401
$ws.Cells.Item(12,7).Value = 1.006277441978455
$ws.Cells.Item(12,8).Value = 1.005635619163513
$ws.Cells.Item(12,9).Value = 1.00449714826917
$ws.Cells.Item(12,10).Value = 1.004676602780819
$ws.Cells.Item(13,7).Value = 1.005270719528198
$ws.Cells.Item(13,8).Value = 1.005754709243774
$ws.Cells.Item(13,9).Value = 1.004501496280567
$ws.Cells.Item(13,10).Value = 1.004674412310123
$ws.Cells.Item(14,7).Value = 1.007476329803467
$ws.Cells.Item(14,8).Value = 1.005756139755249
$ws.Cells.Item(14,9).Value = 1.004945497596367
$ws.Cells.Item(14,10).Value = 1.005081418901682
$ws.Cells.Item(15,7).Value = 1.008578777313232
$ws.Cells.Item(15,8).Value = 1.00580620765686
$ws.Cells.Item(15,9).Value = 1.004984468298129
$ws.Cells.Item(15,10).Value = 1.005116736516356
$ws.Cells.Item(16,7).Value = 1.006074070930481
$ws.Cells.Item(16,8).Value = 1.00568163394928
$ws.Cells.Item(16,9).Value = 1.004470444641429
$ws.Cells.Item(16,10).Value = 1.004504963755608
$ws.Cells.Item(17,7).Value = 1.006881713867188
$ws.Cells.Item(17,8).Value = 1.005273342132568
$ws.Cells.Item(17,9).Value = 1.00445591297171
$ws.Cells.Item(17,10).Value = 1.004562072455883
$ws.Cells.Item(18,7).Value = 1.006578922271729
$ws.Cells.Item(18,8).Value = 1.005491137504578
$ws.Cells.Item(18,9).Value = 1.004466178016695
$ws.Cells.Item(18,10).Value = 1.004518155008554
$ws.Cells.Item(19,7).Value = 1.005913376808167
$ws.Cells.Item(19,8).Value = 1.005580186843872
$ws.Cells.Item(19,9).Value = 1.004458273846915
$ws.Cells.Item(19,10).Value = 1.004541635513306

# Copy formatting (no explicit style) from column G onto columns H:J for rows 2-19
$ws.Range("G2:G19").Copy()
$ws.Range("H2:J19").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---- Caso2 ----
$ws = $wb.Worksheets.Item("Caso2")

# Headers for new AutoML columns (H1, I1, J1)
$ws.Cells.Item(1,8).Value = "Valori stimati autogloun"
$ws.Cells.Item(1,9).Value = "Valori stimati h2o"
$ws.Cells.Item(1,10).Value = "Valori stimati autosklearn"

# Copy formatting (no explicit style) from G1 onto the new header cells
$ws.Cells.Item(1,7).Copy()
$ws.Range($ws.Cells.Item(1,8), $ws.Cells.Item(1,10)).PasteSpecial(-4122)

# Update existing MLP column (G) and populate new AutoML columns (H,I,J) for rows 2-19
$ws.Cells.Item(2,7).Value = 1.009125113487244
$ws.Cells.Item(2,8).Value = 1.004814743995667
$ws.Cells.Item(2,9).Value = 1.00450511169073
$ws.Cells.Item(2,10).Value = 1.004428070038557
$ws.Cells.Item(3,7).Value = 1.005916357040405
$ws.Cells.Item(3,8).Value = 1.004989504814148
$ws.Cells.Item(3,9).Value = 1.004309932200186
$ws.Cells.Item(3,10).Value = 1.004281938076019
$ws.Cells.Item(4,7).Value = 1.006638288497925
$ws.Cells.Item(4,8).Value = 1.005275726318359
$ws.Cells.Item(4,9).Value = 1.004208979579264
$ws.Cells.Item(4,10).Value = 1.004218460991979
$ws.Cells.Item(5,7).Value = 1.005359768867493
$ws.Cells.Item(5,8).Value = 1.005291700363159
$ws.Cells.Item(5,9).Value = 1.004212297416105
$ws.Cells.Item(5,10).Value = 1.004263574257493
$ws.Cells.Item(6,7).Value = 1.005679607391357
$ws.Cells.Item(6,8).Value = 1.005322217941284
$ws.Cells.Item(6,9).Value = 1.00422434839052
$ws.Cells.Item(6,10).Value = 1.004280099645257
$ws.Cells.Item(7,7).Value = 1.004735827445984
$ws.Cells.Item(7,8).Value = 1.005330204963684
$ws.Cells.Item(7,9).Value = 1.004229301495874
$ws.Cells.Item(7,10).Value = 1.004208140075207
$ws.Cells.Item(8,7).Value = 1.005009174346924
$ws.Cells.Item(8,8).Value = 1.005211234092712
$ws.Cells.Item(8,9).Value = 1.004251923094811
$ws.Cells.Item(8,10).Value = 1.004281081259251
$ws.Cells.Item(9,7).Value = 1.005501270294189
$ws.Cells.Item(9,8).Value = 1.005569696426392
$ws.Cells.Item(9,9).Value = 1.00422389712494
$ws.Cells.Item(9,10).Value = 1.004212422296405
$ws.Cells.Item(10,7).Value = 1.005453109741211
$ws.Cells.Item(10,8).Value = 1.005248785018921
$ws.Cells.Item(10,9).Value = 1.00422001624875
$ws.Cells.Item(10,10).Value = 1.004217252135277
$ws.Cells.Item(11,7).Value = 1.004754185676575
$ws.Cells.Item(11,8).Value = 1.005278944969177
$ws.Cells.Item(11,9).Value = 1.004221350352176
$ws.Cells.Item(11,10).Value = 1.004219852387905
$ws.Cells.Item(12,7).Value = 1.005692839622498
$ws.Cells.Item(12,8).Value = 1.005064725875854
$ws.Cells.Item(12,9).Value = 1.004221859468595
$ws.Cells.Item(12,10).Value = 1.004275996237993
$ws.Cells.Item(13,7).Value = 1.004692077636719
$ws.Cells.Item(13,8).Value = 1.005222320556641
$ws.Cells.Item(13,9).Value = 1.004223586657476
$ws.Cells.Item(13,10).Value = 1.004249654710293
$ws.Cells.Item(14,7).Value = 1.006918430328369
$ws.Cells.Item(14,8).Value = 1.005228757858276
$ws.Cells.Item(14,9).Value = 1.004665910024696
$ws.Cells.Item(14,10).Value = 1.004664091393352
$ws.Cells.Item(15,7).Value = 1.008012413978577
$ws.Cells.Item(15,8).Value = 1.005389928817749
$ws.Cells.Item(15,9).Value = 1.00470361998304
$ws.Cells.Item(15,10).Value = 1.004736255854368
$ws.Cells.Item(16,7).Value = 1.005489110946655
$ws.Cells.Item(16,8).Value = 1.005396962165833
$ws.Cells.Item(16,9).Value = 1.004192837746085
$ws.Cells.Item(16,10).Value = 1.00415002182126
$ws.Cells.Item(17,7).Value = 1.006304740905762
$ws.Cells.Item(17,8).Value = 1.005101203918457
$ws.Cells.Item(17,9).Value = 1.004180467405061
$ws.Cells.Item(17,10).Value = 1.004160301759839
$ws.Cells.Item(18,7).Value = 1.005998849868774
$ws.Cells.Item(18,8).Value = 1.005227208137512
$ws.Cells.Item(18,9).Value = 1.004192553226225
$ws.Cells.Item(18,10).Value = 1.004087103530765
$ws.Cells.Item(19,7).Value = 1.005337834358215
$ws.Cells.Item(19,8).Value = 1.005357384681702
$ws.Cells.Item(19,9).Value = 1.004182976158423
$ws.Cells.Item(19,10).Value = 1.004177508875728

# Copy formatting (no explicit style) from column G onto columns H:J for rows 2-19
$ws.Range("G2:G19").Copy()
$ws.Range("H2:J19").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---- Caso3 ----
$ws = $wb.Worksheets.Item("Caso3")

# Headers for new AutoML columns (H1, I1, J1)
$ws.Cells.Item(1,8).Value = "Valori stimati autogloun"
$ws.Cells.Item(1,9).Value = "Valori stimati h2o"
$ws.Cells.Item(1,10).Value = "Valori stimati autosklearn"

# Copy formatting (no explicit style) from G1 onto the new header cells
$ws.Cells.Item(1,7).Copy()
$ws.Range($ws.Cells.Item(1,8), $ws.Cells.Item(1,10)).PasteSpecial(-4122)

# Update existing MLP column (G) and populate new AutoML columns (H,I,J) for rows 2-19
$ws.Cells.Item(2,7).Value = 1.008720874786377
$ws.Cells.Item(2,8).Value = 1.004342555999756
$ws.Cells.Item(2,9).Value = 1.003945836665186
$ws.Cells.Item(2,10).Value = 1.003898121416569
$ws.Cells.Item(3,7).Value = 1.005470156669617
$ws.Cells.Item(3,8).Value = 1.004465222358704
$ws.Cells.Item(3,9).Value = 1.003748131523022
$ws.Cells.Item(3,10).Value = 1.003736581653357
$ws.Cells.Item(4,7).Value = 1.006206512451172
$ws.Cells.Item(4,8).Value = 1.004727840423584
$ws.Cells.Item(4,9).Value = 1.003648794758071
$ws.Cells.Item(4,10).Value = 1.003718765452504
$ws.Cells.Item(5,7).Value = 1.004922747612
$ws.Cells.Item(5,8).Value = 1.004831075668335
$ws.Cells.Item(5,9).Value = 1.003653607168033
$ws.Cells.Item(5,10).Value = 1.003724893555045
$ws.Cells.Item(6,7).Value = 1.005246162414551
$ws.Cells.Item(6,8).Value = 1.00462794303894
$ws.Cells.Item(6,9).Value = 1.003667160928748
$ws.Cells.Item(6,10).Value = 1.003833997994661
$ws.Cells.Item(7,7).Value = 1.00429904460907
$ws.Cells.Item(7,8).Value = 1.004756331443787
$ws.Cells.Item(7,9).Value = 1.003673384641523
$ws.Cells.Item(7,10).Value = 1.003731317818165
$ws.Cells.Item(8,7).Value = 1.004573583602905
$ws.Cells.Item(8,8).Value = 1.004528284072876
$ws.Cells.Item(8,9).Value = 1.003697277007373
$ws.Cells.Item(8,10).Value = 1.003785975277424
$ws.Cells.Item(9,7).Value = 1.005067348480225
$ws.Cells.Item(9,8).Value = 1.004844665527344
$ws.Cells.Item(9,9).Value = 1.003665450876046
$ws.Cells.Item(9,10).Value = 1.003746353089809
$ws.Cells.Item(10,7).Value = 1.005021691322327
$ws.Cells.Item(10,8).Value = 1.004572510719299
$ws.Cells.Item(10,9).Value = 1.003659583574602
$ws.Cells.Item(10,10).Value = 1.003711769357324
$ws.Cells.Item(11,7).Value = 1.00432014465332
$ws.Cells.Item(11,8).Value = 1.004711508750916
$ws.Cells.Item(11,9).Value = 1.00366309137996
$ws.Cells.Item(11,10).Value = 1.003777073696256
$ws.Cells.Item(12,7).Value = 1.005262732505798
$ws.Cells.Item(12,8).Value = 1.004650115966797
$ws.Cells.Item(12,9).Value = 1.003662795665248
$ws.Cells.Item(12,10).Value = 1.003793321549892
$ws.Cells.Item(13,7).Value = 1.004256129264832
$ws.Cells.Item(13,8).Value = 1.004838705062866
$ws.Cells.Item(13,9).Value = 1.003667269187852
$ws.Cells.Item(13,10).Value = 1.003780517727137
$ws.Cells.Item(14,7).Value = 1.006504654884338
$ws.Cells.Item(14,8).Value = 1.004699230194092
$ws.Cells.Item(14,9).Value = 1.004115515334411
$ws.Cells.Item(14,10).Value = 1.004230143502355
$ws.Cells.Item(15,7).Value = 1.007602095603943
$ws.Cells.Item(15,8).Value = 1.004793763160706
$ws.Cells.Item(15,9).Value = 1.004153965173241
$ws.Cells.Item(15,10).Value = 1.004271740093827
$ws.Cells.Item(16,7).Value = 1.00505268573761
$ws.Cells.Item(16,8).Value = 1.004797101020813
$ws.Cells.Item(16,9).Value = 1.003634854582222
$ws.Cells.Item(16,10).Value = 1.003644006326795
$ws.Cells.Item(17,7).Value = 1.005873560905457
$ws.Cells.Item(17,8).Value = 1.004432559013367
$ws.Cells.Item(17,9).Value = 1.003620911222471
$ws.Cells.Item(17,10).Value = 1.0036591719836
$ws.Cells.Item(18,7).Value = 1.005565166473389
$ws.Cells.Item(18,8).Value = 1.004525423049927
$ws.Cells.Item(18,9).Value = 1.003631593950336
$ws.Cells.Item(18,10).Value = 1.003646986559033
$ws.Cells.Item(19,7).Value = 1.004901528358459
$ws.Cells.Item(19,8).Value = 1.004693269729614
$ws.Cells.Item(19,9).Value = 1.003622109830936
$ws.Cells.Item(19,10).Value = 1.003694919869304

# Copy formatting (no explicit style) from column G onto columns H:J for rows 2-19
$ws.Range("G2:G19").Copy()
$ws.Range("H2:J19").PasteSpecial(-4122)
$excel.CutCopyMode = 0

